$p = $ppt.ActivePresentation

function Set-RunText($shape, [string]$text) {
    $tr = $shape.TextFrame.TextRange
    # Force a real text change so the run-merge actually happens even if
    # the final concatenation happens to equal the current text.
    $tr.Text = "TEMP_PLACEHOLDER_TEXT"
    $tr.Text = $text
}

# Slide 1: "Slide 1 (Content)"
Set-RunText $p.Slides.Item(1).Shapes.Item(1) "Slide 1 (Content)"

# Slide 2: "Slide 2 (Content)"
Set-RunText $p.Slides.Item(2).Shapes.Item(1) "Slide 2 (Content)"

# Slide 3: "Slide 3 (Content)"
Set-RunText $p.Slides.Item(3).Shapes.Item(1) "Slide 3 (Content)"

# Slide 4: "Slide 4 (Content)"
Set-RunText $p.Slides.Item(4).Shapes.Item(1) "Slide 4 (Content)"

# Slide 5: "Slide 5 (Two Content)"
Set-RunText $p.Slides.Item(5).Shapes.Item(1) "Slide 5 (Two Content)"

# Slide 6: "Slide 6 (Two Content Right)" and "an image"
Set-RunText $p.Slides.Item(6).Shapes.Item(1) "Slide 6 (Two Content Right)"
Set-RunText $p.Slides.Item(6).Shapes.Item(3) "an image"

# Slide 7: "Slide 7 (Content with Caption)" and "An image"
Set-RunText $p.Slides.Item(7).Shapes.Item(1) "Slide 7 (Content with Caption)"
Set-RunText $p.Slides.Item(7).Shapes.Item(4) "An image"

# Slide 8: "Slide 8 (Comparison)" and "An image"
Set-RunText $p.Slides.Item(8).Shapes.Item(1) "Slide 8 (Comparison)"
Set-RunText $p.Slides.Item(8).Shapes.Item(4) "An image"

# Slide 9: "Slide 9 (Content)"
Set-RunText $p.Slides.Item(9).Shapes.Item(1) "Slide 9 (Content)"

# Slide 10: "Slide 10 (Content)"
Set-RunText $p.Slides.Item(10).Shapes.Item(1) "Slide 10 (Content)"

# Slide 11: "Slide 11 (Content)"
Set-RunText $p.Slides.Item(11).Shapes.Item(1) "Slide 11 (Content)"

# Slide 12: "Slide 12 (Content)"
Set-RunText $p.Slides.Item(12).Shapes.Item(1) "Slide 12 (Content)"
